# fs alert for excel update
#
# The "field type" header (B1) no longer needs the "#client" suffix,
# and the whole "Desc" / description column (E) is being dropped from
# the exported client config sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the B1 header from "字段类型#client" to "字段类型".
$ws.Range("B1").Value = "字段类型"

# 2) Remove column E ("描述"/"Desc"/...) entirely - shifts nothing else,
#    just drops the whole column so the sheet's used range becomes A1:D16.
$ws.Columns("E").Delete()
